# Update the "Kész" (Done) marker for the "Saját Java osztály ..." row
# (row 8) from the "not done" mark (✕) to the "done" mark (✓), matching
# the value chosen from the Munka1!A1:A2 drop-down list used for data
# validation on column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. munkalap - Elkészült funkció")

# Mirror the drop-down selection: Munka1!A2 holds the "done" symbol (✓).
$ws.Range("C8").Value = "✓"

# Move the active selection to C8, as recorded in the saved view state.
$ws.Range("C8").Select()

$wb.Save()
